# Append the new patient record (Suzanne Ismail) as row 5 of the
# "All Patients" sheet, extending the table from A1:H4 to A1:H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = "Suzanne"
$ws.Cells.Item($row, 2).Value = "Ismail"

# surgeryDate is stored as plain text in this sheet (e.g. "2025-05-20"),
# not a real date. Force text formatting before assigning so Excel does
# not auto-convert the literal "2025-07-09" into a date serial number.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "2025-07-09"

$ws.Cells.Item($row, 4).Value = "Lumpectomy"

# notes is blank for this patient, but still needs to exist as an empty
# text cell (like the existing blank notes/goals cells in rows 3-4), so
# use a leading apostrophe to force an explicit empty-text entry instead
# of leaving the cell completely unset.
$ws.Cells.Item($row, 5).Value = "'"

$ws.Cells.Item($row, 6).Value = "Anxiety"

# goals is blank for this patient too.
$ws.Cells.Item($row, 7).Value = "'"

$ws.Cells.Item($row, 8).Value = "suzaism25"
